# Fixed #418 Empty AQL expressions generate empty lines.
# Remove the empty paragraph left behind between "Start of demonstration:"
# and "Some value" (an empty AQL expression used to leave a blank
# paragraph in the generated document; it should be removed entirely).

$d = $word.ActiveDocument

# Locate the empty paragraph (the one whose text is empty/just a
# paragraph mark) and delete its whole range, which removes the
# paragraph (and its paragraph mark) from the document.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "" -or $p.Range.Text -eq "`r") {
        $p.Range.Delete()
    }
}
